$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "67.732.37"
$ws.Range("E2").Value = "  +0.93%  "
Set-TextValue $ws.Range("D3") "2.489.43"
$ws.Range("E3").Value = "  +0.26%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "586.97"
$ws.Range("E5").Value = "  +0.20%  "
Set-TextValue $ws.Range("D6") "176.68"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("E12").Value = "  +0.02%  "
Set-TextValue $ws.Range("D13") "2.942.46"
$ws.Range("E13").Value = "  +0.28%  "
Set-TextValue $ws.Range("D14") "25.66"
$ws.Range("E14").Value = "  +0.94%  "
Set-TextValue $ws.Range("D15") "67.577.46"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  +0.48%  "
Set-TextValue $ws.Range("D17") "2.489.31"
$ws.Range("E17").Value = "  +0.21%  "
Set-TextValue $ws.Range("D18") "7.53"
$ws.Range("E18").Value = "  +1.64%  "
Set-TextValue $ws.Range("D19") "10.94"
$ws.Range("E19").Value = "  -0.43%  "
Set-TextValue $ws.Range("D20") "349.81"
$ws.Range("E20").Value = "  -0.15%  "
Set-TextValue $ws.Range("D21") "4.04"
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("E22").Value = "  -0.05%  "
Set-TextValue $ws.Range("D23") "70.65"
$ws.Range("E23").Value = "  +3.11%  "
Set-TextValue $ws.Range("D24") "4.26"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  -2.72%  "
Set-TextValue $ws.Range("D26") "9.10"
$ws.Range("E26").Value = "  -1.70%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.77%  "
Set-TextValue $ws.Range("D29") "0.0₃0899"
$ws.Range("E29").Value = "  -0.09%  "
Set-TextValue $ws.Range("D30") "505.19"
Set-TextValue $ws.Range("D31") "7.81"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  +2.41%  "
Set-TextValue $ws.Range("D33") "1.77"
$ws.Range("E33").Value = "  +0.27%  "
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  +3.76%  "
Set-TextValue $ws.Range("D36") "162.29"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  -0.11%  "
Set-TextValue $ws.Range("D38") "18.29"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("E42").Value = "  +0.34%  "
Set-TextValue $ws.Range("D43") "4.84"
$ws.Range("E43").Value = "  +0.52%  "
Set-TextValue $ws.Range("D44") "2.41"
$ws.Range("E44").Value = "  +1.10%  "
Set-TextValue $ws.Range("D45") "144.60"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  -0.11%  "
Set-TextValue $ws.Range("D48") "0.0₆0253"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("E50").Value = "  +0.94%  "
Set-TextValue $ws.Range("D51") "0.586"
$ws.Range("E51").Value = "  +0.50%  "
